$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M: this shifts the old M ("Nombre completo acudiente")
# and N ("Teléfono acudiente") columns one place to the right, to N and O
# respectively, leaving the new M column blank (same as Excel's native
# "Insert Sheet Columns" command when column M is selected).
$ws.Columns("M").Insert()

# The freshly inserted column picked up the default column width; give it the
# (~20.4 char / bestFit-less) width Excel used for the newly-inserted blank
# column in the saved workbook.
$ws.Columns("M").ColumnWidth = 19.6

# Put the selection on the newly inserted (now blank) cell the user was
# working with, and scroll the sheet so column E is the left-most visible
# column (matches the recorded view state after the edit).
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("M4").Select() | Out-Null
